$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A86").Value = 85
$ws.Range("B86").Value = 1
$ws.Range("C86").Value = "2024-06-16 17:10:34"
$ws.Range("D86").Value = 200
$ws.Range("E86").Value = 9

$ws.Range("A87").Value = 86
$ws.Range("B87").Value = 2
$ws.Range("C87").Value = "2024-06-16 17:10:34"
$ws.Range("D87").Value = 200
$ws.Range("E87").Value = 1
